# Update Israel MSME country-indicator figures that were refined/recomputed
# with extra decimal precision.
#
#   Enterprises density (per 1000 people): SMEs 4.5 -> 4.54, MSMEs 63 -> 63.04
#   Employment (% of total):     Micro 16.3 -> 16.28, SMEs 36 -> 36.04, MSMEs 52.3 -> 52.32
#   Enterprises (% of total):    Micro 92.4 -> 92.37, SMEs 7.4 -> 7.38, MSMEs 99.7 -> 99.74
#
# These figures are stored as text (not numbers) in the workbook, so the
# target cells are forced to Text format before the new value is written --
# otherwise Excel would auto-coerce the numeric-looking string into a real
# number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "C13" = "4.54"
    "D13" = "63.04"
    "B14" = "16.28"
    "C14" = "36.04"
    "D14" = "52.32"
    "B16" = "92.37"
    "C16" = "7.38"
    "D16" = "99.74"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
}
